# Updated cryptos list on Fri Jun  9 16:31:11 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.502.99"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "'1.843.34"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'262.47"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.5323"
$ws.Range("E7").Value = "  +1.92%  "
$ws.Range("D8").Value = "'0.3108"
$ws.Range("E8").Value = "  -5.12%  "
$ws.Range("D9").Value = "'0.06884"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "'18.57"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").Value = "'0.7626"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "'0.07818"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'1.853.64"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'89.64"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "'5.040"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'14.03"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007942"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "'26.534.54"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "'2.082.98"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").Value = "'4.623"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "'6.014"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'9.306"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").Value = "'141.74"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").Value = "'2.189"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").Value = "'1.687"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "'17.03"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "'110.97"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "'4.283"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").Value = "'0.08782"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "'4.092"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "'0.04831"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "'0.7341"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Value = "'2.929"
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("D36").Value = "'1.134"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").Value = "'3.106"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'2.331"
$ws.Range("E38").Value = "  +5.64%  "
$ws.Range("D39").Value = "'0.01724"
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("D40").Value = "'0.4811"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").Value = "'0.9035"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").Value = "'108.30"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("D43").Value = "'5.899"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'7.645"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "'0.4161"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "'9.026"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "'0.1238"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'34.95"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'0.9002"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05806"
$ws.Range("E51").Value = "  -2.01%  "
